$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number + report week dates) ---
$ws.Range("A8").Characters(21, 1).Text = "8"
$ws.Range("C9").Characters(27, 9).Text = "2/17/2025"
$ws.Range("C9").Characters(47, 9).Text = "2/23/2025"

# --- Weekly CompStat numeric table updates (rows 14-33) ---
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0
$ws.Range("G14").Value = 4
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 5
$ws.Range("J14").Value = 11
$ws.Range("K14").Value = -54.545454545454
$ws.Range("L14").Value = -50
$ws.Range("M14").Value = -66.666666666666
$ws.Range("N14").Value = -92.857142857142
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 50
$ws.Range("F15").Value = 14
$ws.Range("G15").Value = 13
$ws.Range("H15").Value = 7.692307692307
$ws.Range("I15").Value = 34
$ws.Range("J15").Value = 22
$ws.Range("K15").Value = 54.545454545454
$ws.Range("L15").Value = -17.073170731707
$ws.Range("M15").Value = 17.241379310344
$ws.Range("C16").Value = 30
$ws.Range("D16").Value = 47
$ws.Range("E16").Value = -36.170212765957
$ws.Range("F16").Value = 133
$ws.Range("G16").Value = 188
$ws.Range("H16").Value = -29.255319148936
$ws.Range("I16").Value = 230
$ws.Range("J16").Value = 358
$ws.Range("K16").Value = -35.754189944134
$ws.Range("L16").Value = -34.285714285714
$ws.Range("M16").Value = -49.339207048458
$ws.Range("N16").Value = -90.145672664952
$ws.Range("C17").Value = 78
$ws.Range("D17").Value = 61
$ws.Range("E17").Value = 27.868852459016
$ws.Range("F17").Value = 279
$ws.Range("G17").Value = 283
$ws.Range("H17").Value = -1.413427561837
$ws.Range("I17").Value = 517
$ws.Range("J17").Value = 562
$ws.Range("K17").Value = -8.007117437722
$ws.Range("L17").Value = -8.333333333333
$ws.Range("M17").Value = 25.485436893203
$ws.Range("N17").Value = -51.862197392923
$ws.Range("C18").Value = 28
$ws.Range("D18").Value = 31
$ws.Range("E18").Value = -9.677419354838
$ws.Range("F18").Value = 113
$ws.Range("G18").Value = 128
$ws.Range("H18").Value = -11.71875
$ws.Range("I18").Value = 222
$ws.Range("J18").Value = 295
$ws.Range("K18").Value = -24.745762711864
$ws.Range("L18").Value = -34.513274336283
$ws.Range("M18").Value = -42.337662337662
$ws.Range("N18").Value = -86.020151133501
$ws.Range("C19").Value = 91
$ws.Range("D19").Value = 101
$ws.Range("E19").Value = -9.900990099009
$ws.Range("F19").Value = 334
$ws.Range("G19").Value = 415
$ws.Range("H19").Value = -19.518072289156
$ws.Range("I19").Value = 596
$ws.Range("J19").Value = 747
$ws.Range("K19").Value = -20.214190093708
$ws.Range("L19").Value = -31.963470319634
$ws.Range("M19").Value = 18.253968253968
$ws.Range("N19").Value = -28.708133971291
$ws.Range("C20").Value = 23
$ws.Range("D20").Value = 25
$ws.Range("E20").Value = -8
$ws.Range("F20").Value = 85
$ws.Range("G20").Value = 109
$ws.Range("H20").Value = -22.018348623853
$ws.Range("I20").Value = 160
$ws.Range("J20").Value = 236
$ws.Range("K20").Value = -32.203389830508
$ws.Range("L20").Value = -34.426229508196
$ws.Range("M20").Value = -12.087912087912
$ws.Range("N20").Value = -89.382879893828
$ws.Range("C21").Value = 254
$ws.Range("D21").Value = 268
$ws.Range("E21").Value = -5.223880597014
$ws.Range("F21").Value = 962
$ws.Range("G21").Value = 1140
$ws.Range("H21").Value = -15.614035087719
$ws.Range("I21").Value = 1764
$ws.Range("J21").Value = 2231
$ws.Range("K21").Value = -20.932317346481
$ws.Range("L21").Value = -27.227722772277
$ws.Range("M21").Value = -10.95406360424
$ws.Range("N21").Value = -76.461168935148
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 6
$ws.Range("E22").Value = -66.666666666666
$ws.Range("F22").Value = 23
$ws.Range("G22").Value = 28
$ws.Range("H22").Value = -17.857142857142
$ws.Range("I22").Value = 39
$ws.Range("J22").Value = 58
$ws.Range("K22").Value = -32.758620689655
$ws.Range("L22").Value = -22
$ws.Range("M22").Value = -36.065573770491
$ws.Range("C23").Value = 24
$ws.Range("E23").Value = -14.285714285714
$ws.Range("F23").Value = 97
$ws.Range("G23").Value = 100
$ws.Range("H23").Value = -3
$ws.Range("I23").Value = 182
$ws.Range("J23").Value = 205
$ws.Range("K23").Value = -11.219512195122
$ws.Range("L23").Value = -25.409836065573
$ws.Range("M23").Value = 40
$ws.Range("C24").Value = 191
$ws.Range("D24").Value = 240
$ws.Range("E24").Value = -20.416666666666
$ws.Range("F24").Value = 810
$ws.Range("G24").Value = 983
$ws.Range("H24").Value = -17.599186164801
$ws.Range("I24").Value = 1574
$ws.Range("J24").Value = 1733
$ws.Range("K24").Value = -9.174841315637
$ws.Range("L24").Value = -11.173814898419
$ws.Range("M24").Value = 20.336391437308
$ws.Range("C25").Value = 96
$ws.Range("D25").Value = 110
$ws.Range("E25").Value = -12.727272727272
$ws.Range("F25").Value = 328
$ws.Range("G25").Value = 427
$ws.Range("H25").Value = -23.185011709601
$ws.Range("I25").Value = 602
$ws.Range("J25").Value = 730
$ws.Range("K25").Value = -17.534246575342
$ws.Range("L25").Value = -16.272600834492
$ws.Range("C26").Value = 113
$ws.Range("D26").Value = 96
$ws.Range("E26").Value = 17.708333333333
$ws.Range("F26").Value = 420
$ws.Range("G26").Value = 429
$ws.Range("H26").Value = -2.097902097902
$ws.Range("I26").Value = 789
$ws.Range("J26").Value = 832
$ws.Range("K26").Value = -5.168269230769
$ws.Range("L26").Value = -5.95947556615
$ws.Range("M26").Value = -21.1
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 17
$ws.Range("G27").Value = 22
$ws.Range("H27").Value = -22.727272727272
$ws.Range("I27").Value = 44
$ws.Range("J27").Value = 46
$ws.Range("K27").Value = -4.347826086956
$ws.Range("L27").Value = -16.981132075471
$ws.Range("C28").Value = 14
$ws.Range("D28").Value = 10
$ws.Range("E28").Value = 40
$ws.Range("F28").Value = 52
$ws.Range("G28").Value = 41
$ws.Range("H28").Value = 26.829268292682
$ws.Range("I28").Value = 101
$ws.Range("J28").Value = 70
$ws.Range("K28").Value = 44.285714285714
$ws.Range("L28").Value = 21.686746987951
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = -50
$ws.Range("G29").Value = 11
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 23
$ws.Range("J29").Value = 24
$ws.Range("K29").Value = -4.166666666666
$ws.Range("L29").Value = -17.857142857142
$ws.Range("M29").Value = -51.063829787234
$ws.Range("N29").Value = -91.353383458646
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = -50
$ws.Range("G30").Value = 10
$ws.Range("H30").Value = 10
$ws.Range("I30").Value = 20
$ws.Range("J30").Value = 22
$ws.Range("K30").Value = -9.090909090909
$ws.Range("L30").Value = -13.043478260869
$ws.Range("M30").Value = -48.717948717948
$ws.Range("N30").Value = -92.094861660079
$ws.Range("F31").Value = 4
$ws.Range("H31").Value = 300
$ws.Range("L31").Value = -61.538461538461
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = -100
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 2
$ws.Range("J33").Value = 2
